$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.115796208381653
$ws.Range("B1").Value = 2.919003486633301
$ws.Range("C1").Value = 6.640682697296143
$ws.Range("D1").Value = 1.936115503311157
$ws.Range("E1").Value = 1.124441266059875
